$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh
$ws.Cells.Item(2, 4).Value = "26.806.29"
$ws.Cells.Item(2, 5).Value = "  -1.62%  "
$ws.Cells.Item(3, 4).Value = "1.546.20"
$ws.Cells.Item(3, 5).Value = "  -1.82%  "
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "204.53"
$ws.Cells.Item(5, 5).Value = "  -1.39%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.482"
$ws.Cells.Item(6, 5).Value = "  -1.64%  "
$ws.Cells.Item(7, 5).Value = "  +0.14%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "21.39"
$ws.Cells.Item(8, 5).Value = "  -4.34%  "
$ws.Cells.Item(9, 5).Value = "  -1.20%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0580"
$ws.Cells.Item(10, 5).Value = "  -1.96%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0856"
$ws.Cells.Item(11, 5).Value = "  -0.99%  "
$ws.Cells.Item(12, 4).Value = "1.766.48"
$ws.Cells.Item(12, 5).Value = "  -1.83%  "
$ws.Cells.Item(13, 4).Value = "1.548.23"
$ws.Cells.Item(13, 5).Value = "  -1.31%  "
$ws.Cells.Item(14, 5).Value = "  -2.78%  "
$ws.Cells.Item(15, 5).Value = "  -2.18%  "
$ws.Cells.Item(16, 4).Value = "26.786.52"
$ws.Cells.Item(16, 5).Value = "  -1.74%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "60.95"
$ws.Cells.Item(17, 5).Value = "  -2.68%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "213.58"
$ws.Cells.Item(18, 5).Value = "  -1.16%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.25"
$ws.Cells.Item(19, 5).Value = "  -0.70%  "
$ws.Cells.Item(20, 5).Value = "  -1.25%  "
$ws.Cells.Item(21, 5).Value = "  +0.25%  "
$ws.Cells.Item(22, 5).Value = "  -1.58%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.04"
$ws.Cells.Item(23, 5).Value = "  -4.06%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.99"
$ws.Cells.Item(24, 5).Value = "  -0.70%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "152.21"
$ws.Cells.Item(25, 5).Value = "  +0.23%  "
$ws.Cells.Item(26, 5).Value = "  -2.48%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "14.83"
$ws.Cells.Item(27, 5).Value = "  -0.89%  "
$ws.Cells.Item(28, 5).Value = "  +0.10%  "
$ws.Cells.Item(29, 5).Value = "  -2.62%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0461"
$ws.Cells.Item(30, 5).Value = "  -0.72%  "
$ws.Cells.Item(31, 5).Value = "  -3.41%  "
$ws.Cells.Item(32, 5).Value = "  -0.65%  "
$ws.Cells.Item(33, 4).Value = "1.363.59"
$ws.Cells.Item(33, 5).Value = "  -3.01%  "
$ws.Cells.Item(34, 5).Value = "  -0.93%  "
$ws.Cells.Item(35, 5).Value = "  -4.33%  "
$ws.Cells.Item(36, 5).Value = "  -0.60%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.912"
$ws.Cells.Item(37, 5).Value = "  -3.39%  "
$ws.Cells.Item(38, 5).Value = "  -2.30%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.522"
$ws.Cells.Item(39, 5).Value = "  +0.47%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.800"
$ws.Cells.Item(40, 5).Value = "  -2.46%  "
$ws.Cells.Item(41, 5).Value = "  +0.14%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.989"
$ws.Cells.Item(42, 5).Value = "  -1.11%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.49"
$ws.Cells.Item(43, 5).Value = "  +2.57%  "
$ws.Cells.Item(44, 5).Value = "  +0.05%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.77"
$ws.Cells.Item(45, 5).Value = "  -2.59%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "62.70"
$ws.Cells.Item(46, 5).Value = "  -1.94%  "
$ws.Cells.Item(47, 5).Value = "  -1.85%  "
$ws.Cells.Item(48, 4).Value = "1.680.16"
$ws.Cells.Item(48, 5).Value = "  -1.84%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "85.84"
$ws.Cells.Item(49, 5).Value = "  -0.45%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0508"
$ws.Cells.Item(50, 5).Value = "  +2.81%  "
$ws.Cells.Item(51, 4).Value = "0.0₇0960"
$ws.Cells.Item(51, 5).Value = "  -1.46%  "
